# geração de md - novas colunas para série temporal
#
# The original sheet had 6 columns (A:F): ano, total, total_sucesso,
# arrecadado_sucesso, taxa_sucesso, media_sucesso.
# The new layout has 13 columns (A:M):
#   A ano, B total, C total_sucesso, D particip (new), E taxa_sucesso,
#   F arrecadado_sucesso, G media_sucesso, H std_sucesso (new),
#   I min_sucesso (new), J max_sucesso (new), K apoio_medio (new),
#   L contribuicoes (new), M media_contribuicoes (new).
# taxa_sucesso moved from D, now holds e.g. 0.0 - 100.0 (not a 0-1 fraction
# anymore) and keeps a percentage number format; arrecadado_sucesso/
# media_sucesso moved right and keep the currency format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("ano","total","total_sucesso","particip","taxa_sucesso","arrecadado_sucesso","media_sucesso","std_sucesso","min_sucesso","max_sucesso","apoio_medio","contribuicoes","media_contribuicoes")

for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# row data: ano, total, total_sucesso, particip, taxa_sucesso, arrecadado_sucesso,
#           media_sucesso, std_sucesso, min_sucesso, max_sucesso, apoio_medio,
#           contribuicoes, media_contribuicoes
$data = @(
    @(2016, 35, 27, 100, 77.14285714285715, 282854.0616076232, 10476.0763558379, 22534.30873112471, 35.53279454902379, 112920.948828078, 123.8957781899357, 2283, 84.55555555555556),
    @(2017, 85, 61, 100, 71.76470588235294, 547034.1393517209, 8967.77277625772, 15401.01801371438, 34.74344187043801, 95563.29532884162, 84.9695774078473, 6438, 105.5409836065574),
    @(2018, 140, 126, 100, 90, 1220999.723930209, 9690.473999446107, 22024.70522671959, 100.535340323175, 200069.5140664897, 81.37285730957744, 15005, 119.0873015873016),
    @(2019, 159, 148, 100, 93.08176100628931, 1286117.227427096, 8689.981266399294, 15794.77667258089, 29.81192695893366, 152784.2145360522, 86.53146924760114, 14863, 100.4256756756757),
    @(2020, 186, 182, 100, 97.84946236559139, 3201820.20594214, 17592.4187139678, 50496.83999940555, 37.88389448447008, 475290.9541363961, 99.35826860953111, 32225, 177.0604395604396),
    @(2021, 286, 274, 100, 95.8041958041958, 5011848.096713034, 18291.41641136144, 47799.72260745746, 11.93343625774652, 708972.7845446636, 93.12240982372786, 53820, 196.4233576642336),
    @(2022, 308, 305, 100, 99.02597402597402, 4130315.02167561, 13542.0164645102, 30195.27723812832, 10.77163914429046, 374565.1458002281, 82.8764777509804, 49837, 163.4),
    @(2023, 269, 260, 100, 96.6542750929368, 2681143.460911666, 10312.09023427564, 18179.72004853618, 40.21627066051904, 187832.6863604134, 91.8986618992859, 29175, 112.2115384615385)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($col = 1; $col -le $rowVals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowVals[$col - 1]
    }
}

# Number formats: particip (D) is a plain number (the column used to carry
# the currency format pre-edit, so it must be reset back to General);
# taxa_sucesso (E) is a percentage; arrecadado_sucesso (F) and
# media_sucesso (G) are currency (R$). The rest are plain numbers.
$ws.Range("D2:D9").NumberFormat = "General"
$ws.Range("E2:E9").NumberFormat = "0.00%"
$ws.Range("F2:G9").NumberFormat = "R$ #,##0.00"
$ws.Range("H2:M9").NumberFormat = "General"
